# Update fitting parameters on the "Parameters" sheet (r_s_star, h_p_star)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("J2").Value = 0.01005
$ws.Range("K2").Value = 0.241
